$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3882076666666667
$ws.Range("H2").Value = 1.164623
$ws.Range("M2").Value = 0.358199
$ws.Range("N2").Value = 1.074597
$ws.Range("O2").Value = 0.08728684579662628
$ws.Range("P2").Value = 0.08728684579662628
$ws.Range("Q2").Value = 0.1390555979923333
$ws.Range("R2").Value = 1.251500381931
$ws.Range("S2").Value = 0.08728684579662628
$ws.Range("T2").Value = 0.08728684579662628
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3882076666666667
$ws.Range("H3").Value = 1.164623
$ws.Range("O3").Value = 0.004799733963637202
$ws.Range("P3").Value = 0.004799733963637202
$ws.Range("Q3").Value = 0.007646397007777777
$ws.Range("R3").Value = 0.06881757307
$ws.Range("S3").Value = 0.004799733963637202
$ws.Range("T3").Value = 0.004799733963637202
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3882076666666667
$ws.Range("H4").Value = 1.164623
$ws.Range("M4").Value = 0.09277333333333333
$ws.Range("N4").Value = 0.27832
$ws.Range("O4").Value = 0.02260724245658328
$ws.Range("P4").Value = 0.02260724245658328
$ws.Range("Q4").Value = 0.03601531926222222
$ws.Range("R4").Value = 0.32413787336
$ws.Range("S4").Value = 0.02260724245658328
$ws.Range("T4").Value = 0.02260724245658328
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.3882076666666667
$ws.Range("H5").Value = 1.164623
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.3029856666666667
$ws.Range("N5").Value = 0.908957
$ws.Range("O5").Value = 0.07383231992529668
$ws.Range("P5").Value = 0.07383231992529668
$ws.Range("Q5").Value = 0.1176213586901111
$ws.Range("R5").Value = 1.058592228211
$ws.Range("S5").Value = 0.07383231992529668
$ws.Range("T5").Value = 0.07383231992529668
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.3882076666666667
$ws.Range("H6").Value = 1.164623
$ws.Range("M6").Value = 3.282161
$ws.Range("N6").Value = 9.846483000000001
$ws.Range("O6").Value = 0.7998053626244092
$ws.Range("P6").Value = 0.7998053626244092
$ws.Range("Q6").Value = 1.274160063434334
$ws.Range("R6").Value = 11.467440570909
$ws.Range("S6").Value = 0.7998053626244092
$ws.Range("T6").Value = 0.7998053626244092
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.3882076666666667
$ws.Range("H7").Value = 1.164623
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.047884
$ws.Range("N7").Value = 0.143652
$ws.Range("O7").Value = 0.01166849523344748
$ws.Range("P7").Value = 0.01166849523344748
$ws.Range("Q7").Value = 0.01858893591066667
$ws.Range("R7").Value = 0.167300423196
$ws.Range("S7").Value = 0.01166849523344748
$ws.Range("T7").Value = 0.01166849523344748
